$d = $word.ActiveDocument

function Set-ParagraphColorByText($doc, $searchText, $color) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$searchText*") {
            $p.Range.Font.Color = $color
            return $true
        }
    }
    return $false
}

# 1) "Make a simple Currency Converter" loses its red formatting.
Set-ParagraphColorByText $d "Make a simple Currency Converter" ([Microsoft.Office.Interop.Word.WdColor]::wdColorAutomatic)

# 2) "Make an app to display Proximity values" gains red formatting.
Set-ParagraphColorByText $d "Make an app to display Proximity values" ([Microsoft.Office.Interop.Word.WdColor]::wdColorRed)
